# Heroes of Pymoli - Trend Analysis: "Final trend analysis for review"
#
# - Merge the spell-checker-split "Pymoli" runs (w:proofErr spellStart/End
#   wrappers) in the title and in the first bullet back into single plain
#   runs, and extend the first bullet with the new marketing-recommendation
#   sentence.
# - Add a blank "ListParagraph" spacer line after the first bullet.
# - Fill in the text of the second bullet (previously an empty list item)
#   about item purchase price / the Final Critic item.
# - Add two more blank "ListParagraph" spacer lines, then a brand new third
#   bullet about the 20-24 age range, with a couple of Helvetica/18pt styled
#   runs copied out of a spreadsheet ("Average Purchase Price" / "Average
#   Total Purchase Price").

$d = $word.ActiveDocument
$W = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$enDash = [char]0x2013

$blankListParaXml = '<w:p ' + $W + '><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>'

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge the split "Heroes of " / "Pymoli" / " - Trend
#    Analysis" runs into a single plain run.
# ---------------------------------------------------------------------------
$titleXml = '<w:p ' + $W + '><w:r><w:t>Heroes of Pymoli ' + $enDash + ' Trend Analysis</w:t></w:r></w:p>'
$null = $d.Paragraphs(1).Range.InsertXML($titleXml)

# ---------------------------------------------------------------------------
# 2) First bullet: merge the split "...Heroes of " / "Pymoli" / " could ..."
#    runs into one run, and extend the sentence with the new marketing text.
# ---------------------------------------------------------------------------
$bullet1Text = 'Although females only account for approximately 14% of players, they are more inclined to spend more on the purchase of additional items.  Heroes of Pymoli could potentially capture more revenue if they targeted more marketing towards females. '
$bullet1Xml = '<w:p ' + $W + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">' + $bullet1Text + '</w:t></w:r></w:p>'
$null = $d.Paragraphs(2).Range.InsertXML($bullet1Xml)

# ---------------------------------------------------------------------------
# 3) Insert a blank "ListParagraph"-styled paragraph (no numbering) right
#    after the first bullet (paragraph 2).
# ---------------------------------------------------------------------------
$p2End = $d.Paragraphs(2).Range.End
$null = $d.Range($p2End, $p2End).InsertXML($blankListParaXml)

# ---------------------------------------------------------------------------
# 4) The (previously empty) second bullet paragraph - now paragraph 4 - gets
#    its text: two separate runs describing item purchase price / the Final
#    Critic item.
# ---------------------------------------------------------------------------
$run1 = 'The items purchase price does not appear to be a significant deterrent for players.  The '
$run2 = 'Final Critic addition is purchased the most, yet the price is at a higher level.  The items content appears more appealing than price. '
$bullet2Xml = '<w:p ' + $W + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">' + $run1 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $run2 + '</w:t></w:r></w:p>'
$null = $d.Paragraphs(4).Range.InsertXML($bullet2Xml)

# ---------------------------------------------------------------------------
# 5) After that bullet (paragraph 4), insert two blank "ListParagraph"
#    spacer paragraphs followed by a brand-new bullet about the 20-24 age
#    range (mixed plain / Helvetica-styled runs).
# ---------------------------------------------------------------------------
$ageRunIntro = 'The 20-24 age range is the most lucrative in terms of total revenue, however, this age bracket is less inclined to spend on higher priced items, but more inclined on purchasing multiple items.  This is represented by the difference between the '
$helveticaRpr = '<w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="000000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$ageRunAvgPrice = 'Average Purchase Price'
$ageRunAvgTotal = ' and the Average Total Purchase Price'
$ageRunTail = '. '

$ageBulletXml = '<w:p ' + $W + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">' + $ageRunIntro + '</w:t></w:r>' + `
    '<w:r>' + $helveticaRpr + '<w:t>' + $ageRunAvgPrice + '</w:t></w:r>' + `
    '<w:r>' + $helveticaRpr + '<w:t xml:space="preserve">' + $ageRunAvgTotal + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $ageRunTail + '</w:t></w:r></w:p>'

$trailingInsertXml = $blankListParaXml + $blankListParaXml + $ageBulletXml
$p4End = $d.Paragraphs(4).Range.End
$null = $d.Range($p4End, $p4End).InsertXML($trailingInsertXml)

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
